# Disaggregation of commodity "Copper ores and concentrates" -> "Copper".
#
# 1) Rename the commodity label (shared string) from
#    "Copper ores and concentrates" to "Copper".
# 2) For every year sheet (2000..2100), the numeric data in columns D:F
#    (Photovoltaic plants / Offshore wind plants / Onshore wind plants)
#    for rows 5-8 (Neodymium, Dysprosium, Copper, Raw silicon) is
#    re-aligned: each row's three values are cyclically rotated one
#    column to the right, wrapping column F back into column D
#    (new D = old F, new E = old D, new F = old E).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # --- rename the commodity label wherever it appears in column C ---
    for ($r = 5; $r -le 8; $r++) {
        $label = $ws.Cells.Item($r, 3).Value2
        if ($label -eq "Copper ores and concentrates") {
            $ws.Cells.Item($r, 3).Value = "Copper"
        }
    }

    # --- rotate D:F values one column to the right for rows 5-8 ---
    for ($r = 5; $r -le 8; $r++) {
        $d = $ws.Cells.Item($r, 4).Value2
        $e = $ws.Cells.Item($r, 5).Value2
        $f = $ws.Cells.Item($r, 6).Value2

        if ($null -eq $d) { $d = 0 }
        if ($null -eq $e) { $e = 0 }
        if ($null -eq $f) { $f = 0 }

        $ws.Cells.Item($r, 4).Value = $f
        $ws.Cells.Item($r, 5).Value = $d
        $ws.Cells.Item($r, 6).Value = $e
    }
}
